$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.007823705673218
$ws.Range("B1").Value = 2.11767578125
$ws.Range("C1").Value = 6.477560043334961
$ws.Range("D1").Value = 1.680773735046387
$ws.Range("E1").Value = 1.366701126098633
